$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Wnt2"
$ws.Range("C2").Value = "Fzd7"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.01070233333333333
$ws.Range("H2").Value = 0.032107
$ws.Range("I2").Value = 0.006017198313602724
$ws.Range("J2").Value = 0.006017198313602724
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.433028
$ws.Range("N2").Value = 1.299084
$ws.Range("O2").Value = 0.01883770166710565
$ws.Range("P2").Value = 0.01883770166710565
$ws.Range("Q2").Value = 0.004634409998666667
$ws.Range("R2").Value = 0.041709689988
$ws.Range("S2").Value = 0.0001133501867034594
$ws.Range("T2").Value = 0.0001133501867034594

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Wnt2"
$ws.Range("C3").Value = "Fzd7"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.01070233333333333
$ws.Range("H3").Value = 0.032107
$ws.Range("I3").Value = 0.006017198313602724
$ws.Range("J3").Value = 0.006017198313602724
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.993589666666667
$ws.Range("N3").Value = 11.980769
$ws.Range("O3").Value = 0.1737302223447504
$ws.Range("P3").Value = 0.1737302223447504
$ws.Range("Q3").Value = 0.04274072780922222
$ws.Range("R3").Value = 0.384666550283
$ws.Range("S3").Value = 0.001045369200914659
$ws.Range("T3").Value = 0.001045369200914659

# Row 4
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Wnt2"
$ws.Range("C4").Value = "Fzd7"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.01070233333333333
$ws.Range("H4").Value = 0.032107
$ws.Range("I4").Value = 0.006017198313602724
$ws.Range("J4").Value = 0.006017198313602724
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 5.200325666666667
$ws.Range("N4").Value = 15.600977
$ws.Range("O4").Value = 0.2262259795682011
$ws.Range("P4").Value = 0.2262259795682011
$ws.Range("Q4").Value = 0.05565561872655556
$ws.Range("R4").Value = 0.500900568539
$ws.Range("S4").Value = 0.001361246582750904
$ws.Range("T4").Value = 0.001361246582750904

# Row 5
$ws.Range("A5").Value = "ECs"
$ws.Range("B5").Value = "Wnt2"
$ws.Range("C5").Value = "Fzd7"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 1
$ws.Range("F5").Value = 0.3333333333333333
$ws.Range("G5").Value = 0.01070233333333333
$ws.Range("H5").Value = 0.032107
$ws.Range("I5").Value = 0.006017198313602724
$ws.Range("J5").Value = 0.006017198313602724
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 5.344803
$ws.Range("N5").Value = 16.034409
$ws.Range("O5").Value = 0.2325110717631453
$ws.Range("P5").Value = 0.2325110717631453
$ws.Range("Q5").Value = 0.05720186330699999
$ws.Range("R5").Value = 0.514816769763
$ws.Range("S5").Value = 0.00139906522890716
$ws.Range("T5").Value = 0.00139906522890716

# Row 6
$ws.Range("A6").Value = "ECs"
$ws.Range("B6").Value = "Wnt2"
$ws.Range("C6").Value = "Fzd7"
$ws.Range("D6").Value = "Neutro"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.01070233333333333
$ws.Range("H6").Value = 0.032107
$ws.Range("I6").Value = 0.006017198313602724
$ws.Range("J6").Value = 0.006017198313602724
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 3.976918333333334
$ws.Range("N6").Value = 11.930755
$ws.Range("O6").Value = 0.1730049814741227
$ws.Range("P6").Value = 0.1730049814741226
$ws.Range("Q6").Value = 0.04256230564277778
$ws.Range("R6").Value = 0.383060750785
$ws.Range("S6").Value = 0.001041005282770961
$ws.Range("T6").Value = 0.001041005282770961

# Row 7
$ws.Range("A7").Value = "ECs"
$ws.Range("B7").Value = "Wnt2"
$ws.Range("C7").Value = "Fzd7"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.01070233333333333
$ws.Range("H7").Value = 0.032107
$ws.Range("I7").Value = 0.006017198313602724
$ws.Range("J7").Value = 0.006017198313602724
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4.038640666666667
$ws.Range("N7").Value = 12.115922
$ws.Range("O7").Value = 0.1756900431826749
$ws.Range("P7").Value = 0.1756900431826749
$ws.Range("Q7").Value = 0.04322287862822222
$ws.Range("R7").Value = 0.389005907654
$ws.Range("S7").Value = 0.001057161831555581
$ws.Range("T7").Value = 0.001057161831555581

# Row 8
$ws.Range("A8").Value = "FAPs"
$ws.Range("B8").Value = "Wnt2"
$ws.Range("C8").Value = "Fzd7"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 1.767921666666667
$ws.Range("H8").Value = 5.303765
$ws.Range("I8").Value = 0.9939828016863973
$ws.Range("J8").Value = 0.9939828016863973
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 0.433028
$ws.Range("N8").Value = 1.299084
$ws.Range("O8").Value = 0.01883770166710565
$ws.Range("P8").Value = 0.01883770166710565
$ws.Range("Q8").Value = 0.7655595834733334
$ws.Range("R8").Value = 6.890036251260001
$ws.Range("S8").Value = 0.01872435148040219
$ws.Range("T8").Value = 0.01872435148040219

# Row 9
$ws.Range("A9").Value = "FAPs"
$ws.Range("B9").Value = "Wnt2"
$ws.Range("C9").Value = "Fzd7"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 1.767921666666667
$ws.Range("H9").Value = 5.303765
$ws.Range("I9").Value = 0.9939828016863973
$ws.Range("J9").Value = 0.9939828016863973
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.993589666666667
$ws.Range("N9").Value = 11.980769
$ws.Range("O9").Value = 0.1737302223447504
$ws.Range("P9").Value = 0.1737302223447504
$ws.Range("Q9").Value = 7.060353699476112
$ws.Range("R9").Value = 63.543183295285
$ws.Range("S9").Value = 0.1726848531438357
$ws.Range("T9").Value = 0.1726848531438357

# Row 10
$ws.Range("A10").Value = "FAPs"
$ws.Range("B10").Value = "Wnt2"
$ws.Range("C10").Value = "Fzd7"
$ws.Range("D10").Value = "M1"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 1.767921666666667
$ws.Range("H10").Value = 5.303765
$ws.Range("I10").Value = 0.9939828016863973
$ws.Range("J10").Value = 0.9939828016863973
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 5.200325666666667
$ws.Range("N10").Value = 15.600977
$ws.Range("O10").Value = 0.2262259795682011
$ws.Range("P10").Value = 0.2262259795682011
$ws.Range("Q10").Value = 9.19376841982278
$ws.Range("R10").Value = 82.74391577840501
$ws.Range("S10").Value = 0.2248647329854502
$ws.Range("T10").Value = 0.2248647329854502

# Row 11
$ws.Range("A11").Value = "FAPs"
$ws.Range("B11").Value = "Wnt2"
$ws.Range("C11").Value = "Fzd7"
$ws.Range("D11").Value = "M2"
$ws.Range("E11").Value = 3
$ws.Range("F11").Value = 1
$ws.Range("G11").Value = 1.767921666666667
$ws.Range("H11").Value = 5.303765
$ws.Range("I11").Value = 0.9939828016863973
$ws.Range("J11").Value = 0.9939828016863973
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 5.344803
$ws.Range("N11").Value = 16.034409
$ws.Range("O11").Value = 0.2325110717631453
$ws.Range("P11").Value = 0.2325110717631453
$ws.Range("Q11").Value = 9.449193027765
$ws.Range("R11").Value = 85.04273724988501
$ws.Range("S11").Value = 0.2311120065342382
$ws.Range("T11").Value = 0.2311120065342382

# Row 12
$ws.Range("A12").Value = "FAPs"
$ws.Range("B12").Value = "Wnt2"
$ws.Range("C12").Value = "Fzd7"
$ws.Range("D12").Value = "Neutro"
$ws.Range("E12").Value = 3
$ws.Range("F12").Value = 1
$ws.Range("G12").Value = 1.767921666666667
$ws.Range("H12").Value = 5.303765
$ws.Range("I12").Value = 0.9939828016863973
$ws.Range("J12").Value = 0.9939828016863973
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 3.976918333333334
$ws.Range("N12").Value = 11.930755
$ws.Range("O12").Value = 0.1730049814741227
$ws.Range("P12").Value = 0.1730049814741226
$ws.Range("Q12").Value = 7.03088008806389
$ws.Range("R12").Value = 63.27792079257501
$ws.Range("S12").Value = 0.1719639761913517
$ws.Range("T12").Value = 0.1719639761913517

# Row 13
$ws.Range("A13").Value = "FAPs"
$ws.Range("B13").Value = "Wnt2"
$ws.Range("C13").Value = "Fzd7"
$ws.Range("D13").Value = "sCs"
$ws.Range("E13").Value = 3
$ws.Range("F13").Value = 1
$ws.Range("G13").Value = 1.767921666666667
$ws.Range("H13").Value = 5.303765
$ws.Range("I13").Value = 0.9939828016863973
$ws.Range("J13").Value = 0.9939828016863973
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 4.038640666666667
$ws.Range("N13").Value = 12.115922
$ws.Range("O13").Value = 0.1756900431826749
$ws.Range("P13").Value = 0.1756900431826749
$ws.Range("Q13").Value = 7.140000338481112
$ws.Range("R13").Value = 64.26000304633001
$ws.Range("S13").Value = 0.1746328813511193
$ws.Range("T13").Value = 0.1746328813511193
